$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.78905764943711
$ws.Range("C2").Value = 10.55829164071756
$ws.Range("E2").Value = 12.61367878060711
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 28.40749130716947
$ws.Range("H2").Value = 14.35126869846886
$ws.Range("K2").Value = 8.191245009130462
$ws.Range("L2").Value = 9.486424795009569
$ws.Range("N2").Value = 18.6942790829866
$ws.Range("O2").Value = 21.76958335657832
$ws.Range("B3").Value = 11.53588287986293
$ws.Range("C3").Value = 10.58702842652663
$ws.Range("E3").Value = 12.6241952496168
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 28.52728179322251
$ws.Range("H3").Value = 14.39895615079811
$ws.Range("K3").Value = 8.008763643200597
$ws.Range("L3").Value = 9.470397692603806
$ws.Range("N3").Value = 18.74600999576404
$ws.Range("O3").Value = 21.85451155769577
$ws.Range("B4").Value = 11.37940854969903
$ws.Range("C4").Value = 10.60567835477882
$ws.Range("E4").Value = 12.63308338984205
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 28.608761228737
$ws.Range("H4").Value = 14.43018108378491
$ws.Range("K4").Value = 7.895464122005668
$ws.Range("L4").Value = 9.46211421111745
$ws.Range("N4").Value = 18.7794128716341
$ws.Range("O4").Value = 21.91064405734653
$ws.Range("B5").Value = 11.31547652460219
$ws.Range("C5").Value = 10.61353183639095
$ws.Range("E5").Value = 12.63731707673156
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 28.64395136121408
$ws.Range("H5").Value = 14.44339498853765
$ws.Range("K5").Value = 7.849040431043247
$ws.Range("L5").Value = 9.459133031048196
$ws.Range("N5").Value = 18.79343821451918
$ws.Range("O5").Value = 21.93452040855712
$ws.Range("B6").Value = 11.30485323230102
$ws.Range("C6").Value = 10.61485123037962
$ws.Range("E6").Value = 12.63805703011303
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 28.64991443768193
$ws.Range("H6").Value = 14.44561872613269
$ws.Range("K6").Value = 7.841318390682657
$ws.Range("L6").Value = 9.458661906882314
$ws.Range("N6").Value = 18.79579211096286
$ws.Range("O6").Value = 21.93854555740862
$ws.Range("B7").Value = 11.37854689707079
$ws.Range("C7").Value = 10.60578324234413
$ws.Range("E7").Value = 12.63313800969144
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 28.60922778069388
$ws.Range("H7").Value = 14.43035730853105
$ws.Range("K7").Value = 7.894838976869849
$ws.Range("L7").Value = 9.462072405405769
$ws.Range("N7").Value = 18.77960034709158
$ws.Range("O7").Value = 21.91096200611712
$ws.Range("B8").Value = 11.70202860747163
$ws.Range("C8").Value = 10.56799180709182
$ws.Range("E8").Value = 12.61680061007226
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 28.44714562994761
$ws.Range("H8").Value = 14.36730801809504
$ws.Range("K8").Value = 8.128622622024505
$ws.Range("L8").Value = 9.480577223876868
$ws.Range("N8").Value = 18.71177616862808
$ws.Range("O8").Value = 21.79803897366853
$ws.Range("B9").Value = 12.32445706161763
$ws.Range("C9").Value = 10.50183041386077
$ws.Range("E9").Value = 12.60402663818918
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 28.19250337593389
$ws.Range("H9").Value = 14.25907604800665
$ws.Range("K9").Value = 8.574504339912071
$ws.Range("L9").Value = 9.529085878854625
$ws.Range("N9").Value = 18.5917372526017
$ws.Range("O9").Value = 21.60824718399774
$ws.Range("B10").Value = 12.76967693100752
$ws.Range("C10").Value = 10.45802465453402
$ws.Range("E10").Value = 12.60634046502297
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 28.0443418757775
$ws.Range("H10").Value = 14.18891977496938
$ws.Range("K10").Value = 8.891151576466024
$ws.Range("L10").Value = 9.571978622516404
$ws.Range("N10").Value = 18.51138070215102
$ws.Range("O10").Value = 21.48813018075391
$ws.Range("B11").Value = 12.96866399138912
$ws.Range("C11").Value = 10.43913032445834
$ws.Range("E11").Value = 12.60991934021326
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 27.98547133181585
$ws.Range("H11").Value = 14.15903030607067
$ws.Range("K11").Value = 9.032197960916269
$ws.Range("L11").Value = 9.593019343002997
$ws.Range("N11").Value = 18.47651173861368
$ws.Range("O11").Value = 21.4376880723641
$ws.Range("B12").Value = 13.04343403099317
$ws.Range("C12").Value = 10.43212342845538
$ws.Range("E12").Value = 12.61163624514774
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 27.96441121752696
$ws.Range("H12").Value = 14.14800269666661
$ws.Range("K12").Value = 9.085129717031757
$ws.Range("L12").Value = 9.601202295476597
$ws.Range("N12").Value = 18.46354912113785
$ws.Range("O12").Value = 21.41919151323565
$ws.Range("B13").Value = 13.02735794820072
$ws.Range("C13").Value = 10.4336259168373
$ws.Range("E13").Value = 12.61125042218432
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 27.96889195783295
$ws.Range("H13").Value = 14.15036475723906
$ws.Range("K13").Value = 9.073751971711406
$ws.Range("L13").Value = 9.599430449852688
$ws.Range("N13").Value = 18.46633012588632
$ws.Range("O13").Value = 21.42314816701254
$ws.Range("B14").Value = 12.97482742835409
$ws.Range("C14").Value = 10.43855090056599
$ws.Range("E14").Value = 12.61005335286323
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 27.98371395891011
$ws.Range("H14").Value = 14.15811722996228
$ws.Range("K14").Value = 9.03656256143582
$ws.Range("L14").Value = 9.593688271502335
$ws.Range("N14").Value = 18.47544046125477
$ws.Range("O14").Value = 21.43615422382242
$ws.Range("B15").Value = 12.94257311150344
$ws.Range("C15").Value = 10.44158684896608
$ws.Range("E15").Value = 12.60936716044655
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 27.99295360631278
$ws.Range("H15").Value = 14.1629037161905
$ws.Range("K15").Value = 9.013719161617399
$ws.Range("L15").Value = 9.590198919951447
$ws.Range("N15").Value = 18.48105222927533
$ws.Range("O15").Value = 21.44419959402199
$ws.Range("B16").Value = 12.7565950278155
$ws.Range("C16").Value = 10.45928018240049
$ws.Range("E16").Value = 12.60615728001281
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 28.04836138448792
$ws.Range("H16").Value = 14.1909138433165
$ws.Range("K16").Value = 8.881869322476355
$ws.Range("L16").Value = 9.570633940307637
$ws.Range("N16").Value = 18.51369331380212
$ws.Range("O16").Value = 21.49151127339627
$ws.Range("B17").Value = 12.64154244388737
$ws.Range("C17").Value = 10.47039865955187
$ws.Range("E17").Value = 12.60483409428665
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 28.08454138793433
$ws.Range("H17").Value = 14.20861557191046
$ws.Range("K17").Value = 8.8001803301704
$ws.Range("L17").Value = 9.55901983070572
$ws.Range("N17").Value = 18.53414865566373
$ws.Range("O17").Value = 21.5216117126798
$ws.Range("B18").Value = 12.57503750884102
$ws.Range("C18").Value = 10.47689099120718
$ws.Range("E18").Value = 12.60431093231264
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 28.10615357936029
$ws.Range("H18").Value = 14.2189877552895
$ws.Range("K18").Value = 8.752915458122327
$ws.Range("L18").Value = 9.552483867391846
$ws.Range("N18").Value = 18.54607273946461
$ws.Range("O18").Value = 21.53931993471378
$ws.Range("B19").Value = 12.55246572878155
$ws.Range("C19").Value = 10.47910590864837
$ws.Range("E19").Value = 12.60417470349118
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 28.11360870956384
$ws.Range("H19").Value = 14.22253234656867
$ws.Range("K19").Value = 8.73686590702707
$ws.Range("L19").Value = 9.550295796773085
$ws.Range("N19").Value = 18.55013731567616
$ws.Range("O19").Value = 21.54538349222822
$ws.Range("B20").Value = 12.65382465234845
$ws.Range("C20").Value = 10.46920501507312
$ws.Range("E20").Value = 12.60495033892959
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 28.08060687071223
$ws.Range("H20").Value = 14.20671146675558
$ws.Range("K20").Value = 8.808905553969572
$ws.Range("L20").Value = 9.560241283903165
$ws.Range("N20").Value = 18.53195472997541
$ws.Range("O20").Value = 21.51836655691856
$ws.Range("B21").Value = 12.99027326250536
$ws.Range("C21").Value = 10.43710030205363
$ws.Range("E21").Value = 12.61039515978439
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 27.97932686798544
$ws.Range("H21").Value = 14.15583224866138
$ws.Range("K21").Value = 9.04749937239834
$ws.Range("L21").Value = 9.595369082867231
$ws.Range("N21").Value = 18.47275798763984
$ws.Range("O21").Value = 21.43231760942979
$ws.Range("B22").Value = 13.20673931264
$ws.Range("C22").Value = 10.41698027059244
$ws.Range("E22").Value = 12.61606079353401
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 27.92032348582237
$ws.Range("H22").Value = 14.12427509962459
$ws.Range("K22").Value = 9.200619713977511
$ws.Range("L22").Value = 9.619579726821383
$ws.Range("N22").Value = 18.43547675146809
$ws.Range("O22").Value = 21.37960512944989
$ws.Range("B23").Value = 13.09154303772231
$ws.Range("C23").Value = 10.42764000561199
$ws.Range("E23").Value = 12.61284473125791
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 27.95115486309703
$ws.Range("H23").Value = 14.14096271556632
$ws.Range("K23").Value = 9.119169085189293
$ws.Range("L23").Value = 9.606544994483542
$ws.Range("N23").Value = 18.45524597512647
$ws.Range("O23").Value = 21.40741590454564
$ws.Range("B24").Value = 12.64827298065627
$ws.Range("C24").Value = 10.46974434967953
$ws.Range("E24").Value = 12.60489704465835
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 28.08238313767583
$ws.Range("H24").Value = 14.20757170456407
$ws.Range("K24").Value = 8.804961813444301
$ws.Range("L24").Value = 9.559688624933026
$ws.Range("N24").Value = 18.53294609275094
$ws.Range("O24").Value = 21.51983243634385
$ws.Range("B25").Value = 12.15786284189264
$ws.Range("C25").Value = 10.5188823548199
$ws.Range("E25").Value = 12.6054237626461
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 28.25458305792943
$ws.Range("H25").Value = 14.28670944704898
$ws.Range("K25").Value = 8.455578197163934
$ws.Range("L25").Value = 9.514673919812648
$ws.Range("N25").Value = 18.62283006471985
$ws.Range("O25").Value = 21.65619966883991
